$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New shared string for the extra "Contact #" column header. Move the
# existing "Emails" header from C5 to the new D5 first so the shared string
# for "Emails" is preserved, then overwrite C5 with the new header text.
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Emails"
$ws.Range("C5").Value = "Contact #"

# ---------------------------------------------------------------------------
# Column C (Contacts numbers) is formatted as Text, left aligned.
# ---------------------------------------------------------------------------
$ws.Range("C2:C5").NumberFormat = "@"
$ws.Range("C2:C5").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Columns A, B and D (besides the title/header rows) are left+top aligned.
# ---------------------------------------------------------------------------
$ws.Range("A2:B4").HorizontalAlignment = -4131
$ws.Range("A2:B4").VerticalAlignment = -4160
$ws.Range("D1:D4").HorizontalAlignment = -4131
$ws.Range("D1:D4").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Header row (row 5) - bold, boxed with a thin border on every side.
# ---------------------------------------------------------------------------
$ws.Range("A5:D5").Font.Bold = $true
$ws.Range("A5:D5").Borders.LineStyle = 1
$ws.Range("A5:B5").HorizontalAlignment = -4131
$ws.Range("A5:B5").VerticalAlignment = -4160
$ws.Range("D5").HorizontalAlignment = -4131
$ws.Range("D5").VerticalAlignment = -4160
$ws.Range("C5").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# New column D is as wide as the export needs for an email address.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 56

# ---------------------------------------------------------------------------
# Dimension grows to D5, and the active selection moves to column E (the
# next free column) the way Excel leaves it after inserting a column.
# ---------------------------------------------------------------------------
$ws.Range("E1:E1048576").Select()

Write-Host "done"
